$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save" — copy the formatting used by the other
# header cells (bold/border/centered style) from G1.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data column values
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
